$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "Partial visibility" results table (rows 29-41) ---
# The grid-size label used for this run of the partially-observable
# simulation changed from 5x5 to 15x15 ...
$ws.Range("Z29").Value = "grid size 15 by 15"
# ... and the final (10th) run's grid-size label swapped the other way,
# 15x15 -> 5x5.
$ws.Range("Z41").Value = "grid size 5 by 5"

# Fill in the previously-empty "pits hit / bonuses collected / deaths"
# results (columns AA:AC) for runs 1-10 of the partial-visibility table,
# now that those simulations have been run.
$results = @(
    @(5,  12, 1),
    @(2,   1, 1),
    @(59, 39, 0),
    @(12, 16, 1),
    @(2,   1, 1),
    @(25, 22, 1),
    @(14,  7, 0),
    @(5,   2, 0),
    @(53, 46, 2),
    @(97, 58, 2)
)

$startRow = 31
for ($i = 0; $i -lt $results.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("AA$row").Value = $results[$i][0]
    $ws.Range("AB$row").Value = $results[$i][1]
    $ws.Range("AC$row").Value = $results[$i][2]
}

# Update the saved scroll position / active selection of the sheet to
# reflect where the author was last working.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 16
$ws.Range("AC40").Select() | Out-Null
